$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (matches source data which is not numeric)
$textCells = @('D5', 'D6', 'D8', 'D9', 'D10', 'D11', 'D12', 'D15', 'D17', 'D20', 'D21', 'D22', 'D23', 'D24', 'D25', 'D27', 'D28', 'D30', 'D31', 'D33', 'D36', 'D38', 'D39', 'D40', 'D42', 'D44', 'D46', 'D48', 'D49', 'D50')
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '26.141.19'
$ws.Range('E2').Value = '  -0.51%  '

$ws.Range('D3').Value = '1.655.59'
$ws.Range('E3').Value = '  -0.70%  '

$ws.Range('E4').Value = '  -0.42%  '

$ws.Range('D5').Value = '219.37'
$ws.Range('E5').Value = '  -0.22%  '

$ws.Range('D6').Value = '0.5261'

$ws.Range('E7').Value = '  -0.41%  '

$ws.Range('D8').Value = '0.2688'
$ws.Range('E8').Value = '  +1.44%  '

$ws.Range('D9').Value = '0.06379'
$ws.Range('E9').Value = '  +0.10%  '

$ws.Range('D10').Value = '20.59'
$ws.Range('E10').Value = '  -1.64%  '

$ws.Range('D11').Value = '0.07700'
$ws.Range('E11').Value = '  -1.76%  '

$ws.Range('D12').Value = '4.616'
$ws.Range('E12').Value = '  +2.03%  '

$ws.Range('D13').Value = '1.663.87'
$ws.Range('E13').Value = '  -0.38%  '

$ws.Range('D14').Value = '1.884.12'

$ws.Range('D15').Value = '0.5640'
$ws.Range('E15').Value = '  +0.72%  '

$ws.Range('D16').Value = '0.0₅8260'
$ws.Range('E16').Value = '  +1.80%  '

$ws.Range('D17').Value = '65.73'

$ws.Range('D18').Value = '26.124.94'
$ws.Range('E18').Value = '  -0.61%  '

$ws.Range('E19').Value = '  -0.43%  '

$ws.Range('D20').Value = '4.693'
$ws.Range('E20').Value = '  -0.82%  '

$ws.Range('D21').Value = '10.37'
$ws.Range('E21').Value = '  +1.08%  '

$ws.Range('D22').Value = '190.50'
$ws.Range('E22').Value = '  -4.84%  '

$ws.Range('D23').Value = '6.006'
$ws.Range('E23').Value = '  -1.09%  '

$ws.Range('D24').Value = '1.004'
$ws.Range('E24').Value = '  -0.50%  '

$ws.Range('D25').Value = '151.25'
$ws.Range('E25').Value = '  +3.66%  '

$ws.Range('E26').Value = '  -1.02%  '

$ws.Range('D27').Value = '7.277'
$ws.Range('E27').Value = '  +0.56%  '

$ws.Range('D28').Value = '16.03'
$ws.Range('E28').Value = '  -1.22%  '

$ws.Range('E29').Value = '  -0.55%  '

$ws.Range('D30').Value = '0.05654'
$ws.Range('E30').Value = '  -4.23%  '

$ws.Range('D31').Value = '1.278'
$ws.Range('E31').Value = '  -0.32%  '

$ws.Range('E32').Value = '  -0.37%  '

$ws.Range('D33').Value = '3.388'
$ws.Range('E33').Value = '  +1.94%  '

$ws.Range('E34').Value = '  -0.85%  '

$ws.Range('E35').Value = '  -0.74%  '

$ws.Range('D36').Value = '0.9496'
$ws.Range('E36').Value = '  -1.20%  '

$ws.Range('E37').Value = '  -0.79%  '

$ws.Range('D38').Value = '0.5785'
$ws.Range('E38').Value = '  -0.20%  '

$ws.Range('D39').Value = '0.01600'
$ws.Range('E39').Value = '  -0.83%  '

$ws.Range('D40').Value = '5.977'
$ws.Range('E40').Value = '  +0.14%  '

$ws.Range('E41').Value = '  -0.44%  '

$ws.Range('D42').Value = '0.8349'
$ws.Range('E42').Value = '  -2.72%  '

$ws.Range('D43').Value = '1.025.43'
$ws.Range('E43').Value = '  -4.64%  '

$ws.Range('D44').Value = '101.47'
$ws.Range('E44').Value = '  -1.29%  '

$ws.Range('D45').Value = '1.793.93'
$ws.Range('E45').Value = '  -0.66%  '

$ws.Range('D46').Value = '58.45'
$ws.Range('E46').Value = '  -0.09%  '

$ws.Range('E47').Value = '  +1.98%  '

$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = '0.05346'
$ws.Range('E48').Value = '  +3.98%  '

$ws.Range('B49').Value = 'Frax'
$ws.Range('C49').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D49').Value = '1.004'
$ws.Range('E49').Value = '  -1.04%  '

$ws.Range('D50').Value = '8.058'
$ws.Range('E50').Value = '  -0.48%  '

$ws.Range('E51').Value = '  -1.58%  '

foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).Style = "Normal"
}
